$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of portfolio data appended at row 33.
# The date must remain a plain text string (matching the existing rows),
# so force text entry then clear the resulting formatting/style so the
# new cells match the unstyled look of the other data rows.
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "2025-09-17"
$ws.Range("A33").ClearFormats()

$ws.Range("B33").Value = 59.29000091552734
$ws.Range("C33").Value = 719.1500244140625
$ws.Range("D33").Value = 328.25
